$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the whole row 34 (shifts rows 35:113 up to 34:112, and adjusts formulas)
$ws.Rows("34").Delete()

# Update the view to match the target state
$excel.ActiveWindow.ScrollRow = 31
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A34:XFD34").Select()
